$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Proveedor (column B) values
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1

# Update Puntos (column E) values
$ws.Range("E2").Value = 1250
$ws.Range("E3").Value = 1300
$ws.Range("E4").Value = 1399

# Update the active selection on the sheet view
$ws.Range("A5").Select()
